$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 43123
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "weekly meeting and playing around with data in R"

$ws.Range("A7").Select()
